$d = $word.ActiveDocument

# 1. Rename the title from "Milestone 4" to "Milestone 6"
$d.Content.Find.Execute("Milestone 4", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Milestone 6", 2)

# 2. Remove the existing "_GoBack" bookmark (currently sitting at the end of
#    the "Buying houses, jail, taxes, and win/loss." paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3. Re-add the "_GoBack" bookmark around the title run "Milestone 6" in the
#    first paragraph (matching where the bookmark now lives in the diff).
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
[void]$r.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $r)
